$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the formatting
# (bold font, border, centered alignment) from the existing header H1.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF) for rows 2-27
$data = @{
    2  = @(8, 9)
    3  = @(1, 4)
    4  = @(1, 8)
    5  = @(1, 4)
    6  = @(7, 7)
    7  = @(7, 9)
    8  = @(5, 7)
    9  = @(1, 3)
    10 = @(1, 4)
    11 = @(1, 5)
    12 = @(1, 5)
    13 = @(1, 6)
    14 = @(1, 5)
    15 = @(1, 6)
    16 = @(1, 7)
    17 = @(1, 6)
    18 = @(1, 7)
    19 = @(1, 5)
    20 = @(1, 7)
    21 = @(1, 6)
    22 = @(1, 4)
    23 = @(1, 3)
    24 = @(1, 4)
    25 = @(1, 4)
    26 = @(1, 3)
    27 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
